# Apply the "Adding US & Metric variants" edit:
# - Swap the "Scrape" value for "Tank Cleaning Nozzles" (B4) and "Air Nozzles" (B6)
# - Move active cell selection from B7 to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Swap B4 and B6 values ("No" <-> "Yes")
$b4 = $ws.Range("B4").Value2
$b6 = $ws.Range("B6").Value2
$ws.Range("B4").Value = $b6
$ws.Range("B6").Value = $b4

# Update the active selection to B4
$ws.Activate()
$ws.Range("B4").Select()

# Note: the workbook window geometry (xWindow/yWindow/windowWidth/windowHeight)
# and the workbook's internal revisionPtr documentId GUID reflect the host
# Excel session's screen/window state at save time. They are not exposed
# through the scriptable object model (confirmed empirically: setting
# Window.Left/Top/Width/Height, WindowState, or Application.Left/Top/Width/
# Height has no effect on the serialized workbookView element in this
# runtime), so they are intentionally left alone here.
